$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# Slide 1 : Cover (Title + Subtitle)
# ----------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Sanyam Purohit: MCA Student Profile"
$s1.Shapes.Item(2).TextFrame.TextRange.Runs(1, 1).Text = "Source: Sanyam_Purohit_Resume"

# ----------------------------------------------------------------------
# Slide 2 : Contents list (7 bullets -> 6 bullets)
# ----------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Paragraphs(2, 1).Runs(1, 1).Text = "Contact Information"
$tr2.Paragraphs(3, 1).Runs(1, 1).Text = "Summary"
$tr2.Paragraphs(4, 1).Runs(1, 1).Text = "Education"
$tr2.Paragraphs(5, 1).Runs(1, 1).Text = "Technical Skills"
$tr2.Paragraphs(6, 1).Runs(1, 1).Text = "Projects"
$tr2.Paragraphs(7, 1).Runs(1, 1).Text = "Certifications & Workshops"
$tr2.Paragraphs(8, 1).Delete()

# ----------------------------------------------------------------------
# Slide 3 : Section header -> Contact Information
# ----------------------------------------------------------------------
$p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Contact Information"

# ----------------------------------------------------------------------
# Slide 4 : Contact Information detail (5 bullets -> 3 bullets)
# ----------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Contact Information"
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Paragraphs(2, 1).Runs(1, 1).Text = "Sabarkantha, Gujarat"
$tr4.Paragraphs(3, 1).Runs(1, 1).Text = "purohitsanyam0311@gmail.com"
$tr4.Paragraphs(4, 1).Runs(1, 1).Text = "+91 9408814497"
$tr4.Paragraphs(6, 1).Delete()
$tr4.Paragraphs(5, 1).Delete()

# ----------------------------------------------------------------------
# Slide 5 : Section header -> Summary
# ----------------------------------------------------------------------
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Summary"

# ----------------------------------------------------------------------
# Slide 6 : Summary detail (5 bullets -> 4 bullets)
# ----------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Summary"
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange
$tr6.Paragraphs(2, 1).Runs(1, 1).Text = "First-year MCA student focused on software development"
$tr6.Paragraphs(3, 1).Runs(1, 1).Text = "Strengths in programming, problem-solving, and learning new technologies"
$tr6.Paragraphs(4, 1).Runs(1, 1).Text = "Passionate about hands-on projects and real-world applications"
$tr6.Paragraphs(5, 1).Runs(1, 1).Text = "Actively seeking growth and contribution opportunities"
$tr6.Paragraphs(6, 1).Delete()

# ----------------------------------------------------------------------
# Slide 7 : Section header -> Education
# ----------------------------------------------------------------------
$p.Slides.Item(7).Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Education"

# ----------------------------------------------------------------------
# Slide 8 : Education detail (5 bullets -> 3 bullets)
# ----------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Education"
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange
$tr8.Paragraphs(2, 1).Runs(1, 1).Text = "Master of Computer Applications (MCA): 2024 - Current, Kadi Sarva Vishwavidyalaya (Gujarat)"
$tr8.Paragraphs(3, 1).Runs(1, 1).Text = "Bachelor of Computer Applications (BCA): 2021 - 2024, HNGU (Gujarat), CGPA: 9.28"
$tr8.Paragraphs(4, 1).Runs(1, 1).Text = "Higher Secondary (XII): 2020 - 2021, GHSEB (Gujarat), Percentage: 83.57%"
$tr8.Paragraphs(6, 1).Delete()
$tr8.Paragraphs(5, 1).Delete()

# ----------------------------------------------------------------------
# Slide 9 : Section header -> Technical Skills
# ----------------------------------------------------------------------
$p.Slides.Item(9).Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Technical Skills"

# ----------------------------------------------------------------------
# Slide 10 : Technical Skills detail (5 bullets -> 4 bullets)
# ----------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Technical Skills"
$tr10 = $s10.Shapes.Item(2).TextFrame.TextRange
$tr10.Paragraphs(2, 1).Runs(1, 1).Text = "Frontend: HTML, CSS"
$tr10.Paragraphs(3, 1).Runs(1, 1).Text = "Backend: PHP, Java, Python, ASP.NET with C#, Swift"
$tr10.Paragraphs(4, 1).Runs(1, 1).Text = "Database: MySql, Oracle, SQL Server"
$tr10.Paragraphs(5, 1).Runs(1, 1).Text = "Tools: Visual Studio, Eclipse, XAMPP, Xcode, Android Studio"
$tr10.Paragraphs(6, 1).Delete()

# ----------------------------------------------------------------------
# Slide 11 : Section header -> Projects
# ----------------------------------------------------------------------
$p.Slides.Item(11).Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Projects"

# ----------------------------------------------------------------------
# Slide 12 : Projects detail (4 bullets -> 7 bullets)
# ----------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Projects"
$tr12 = $s12.Shapes.Item(2).TextFrame.TextRange
$tr12.Paragraphs(2, 1).Runs(1, 1).Text = "VASUDHA (Online Society Management System)"
$tr12.Paragraphs(3, 1).Runs(1, 1).Text = "ASP.NET (C#) | 2.5 months"
$tr12.Paragraphs(4, 1).Runs(1, 1).Text = "Streamlines society operations and member-management communication"
$tr12.Paragraphs(5, 1).Runs(1, 1).Text = "GOODKARMA (Micro Donation Platform)"
$lastPara12 = $tr12.Paragraphs($tr12.Paragraphs().Count, 1)
$lastPara12.InsertAfter("`rASP.NET (C#) | 3 months`rFeatures: User authentication, donation tracking, admin/donor roles`rSupports transparent micro-donations for social causes")

# ----------------------------------------------------------------------
# Slide 13 : Section header -> Certifications & Workshops
# ----------------------------------------------------------------------
$p.Slides.Item(13).Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Certifications & Workshops"

# ----------------------------------------------------------------------
# Slide 14 : Certifications & Workshops detail (3 bullets -> 2 bullets)
# ----------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$s14.Shapes.Item(1).TextFrame.TextRange.Runs(1, 1).Text = "Certifications & Workshops"
$tr14 = $s14.Shapes.Item(2).TextFrame.TextRange
$tr14.Paragraphs(2, 1).Runs(1, 1).Text = "5-Day Workshop on Laravel & WordPress: Jan 2025 (LDRP Institute, Prof. Adarsh Patel)"
$tr14.Paragraphs(3, 1).Runs(1, 1).Text = "AI for Students – Build Your Own Generative AI Model: Sep 2023 (NXT WAVE/IIT Delhi)"
$tr14.Paragraphs(4, 1).Delete()

# ----------------------------------------------------------------------
# Slide 15 : Section header -> Conclusion (unchanged)
# Slide 16 : Conclusion detail - reword bullets & drop trailing empty paragraph
# ----------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$tr16 = $s16.Shapes.Item(2).TextFrame.TextRange
$tr16.Paragraphs(2, 1).Runs(1, 1).Text = "Aspiring developer with strong academic foundation and hands-on project experience."
$tr16.Paragraphs(3, 1).Runs(1, 1).Text = "Demonstrates proficiency in diverse technologies and eagerness for continuous growth."
$tr16.Paragraphs(4, 1).Delete()
